$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Raluca"
$ws.Range("B2").Value = "Sofrone"
$ws.Range("D2").Value = "Team Leader Marketing"
$ws.Range("E2").Value = "Quartz Matrix"

$ws.Columns.Item(1).ColumnWidth = 11.67
$ws.Columns.Item(2).ColumnWidth = 12.67

$ws.Range("D9").Select()
